# "Test Cases" sheet (xl/worksheets/sheet1.xml) gets a new test-case row
# (AuthoringRecordViewDetailsTest) plus a couple of result/runmode tweaks on
# the existing rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Row 2 (AuthoringTest): Results SKIP -> FAIL
$ws.Range("D2").Value = "FAIL"

# Runmode Y -> N for the three existing scenario rows
$ws.Range("C3").Value = "N"
$ws.Range("C4").Value = "N"
$ws.Range("C5").Value = "N"

# New row 6: AuthoringRecordViewDetailsTest
$ws.Range("A6").Value = "AuthoringRecordViewDetailsTest"
$ws.Range("B6").Value = "To verify Record View Details link Navigate to WOS page and Navigate to Project Neon Page"
$ws.Range("C6").Value = "N"
$ws.Range("D6").Value = "SKIP"

# Match the new row's look-and-feel (borders) to the rest of the table
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("D5").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("D5").Copy()
$ws.Range("D6").PasteSpecial(-4122)

# Move the sheet's active selection to C12
$ws.Range("C12").Select() | Out-Null
